$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '45.409.40'; E = '  +6.88%  ' },
    @{ Row = 3; D = '2.375.86'; E = '  +3.92%  ' },
    @{ Row = 4; E = '  +0.24%  ' },
    @{ Row = 5; D = '111.53'; E = '  +7.75%  ' },
    @{ Row = 6; D = '317.73'; E = '  +2.12%  ' },
    @{ Row = 7; E = '  +2.46%  ' },
    @{ Row = 8; E = '  -0.13%  ' },
    @{ Row = 9; D = '0.631'; E = '  +5.37%  ' },
    @{ Row = 10; D = '41.92'; E = '  +7.64%  ' },
    @{ Row = 11; D = '0.0930'; E = '  +3.18%  ' },
    @{ Row = 12; E = '  +5.23%  ' },
    @{ Row = 13; E = '  +3.88%  ' },
    @{ Row = 14; E = '  +0.40%  ' },
    @{ Row = 15; D = '15.79'; E = '  +4.81%  ' },
    @{ Row = 16; D = '2.739.87'; E = '  +3.82%  ' },
    @{ Row = 17; D = '2.370.80'; E = '  +3.73%  ' },
    @{ Row = 18; D = '45.222.86'; E = '  +6.02%  ' },
    @{ Row = 19; E = '  +5.00%  ' },
    @{ Row = 20; E = '  +3.65%  ' },
    @{ Row = 21; D = '13.03'; E = '  -4.06%  ' },
    @{ Row = 22; D = '75.14'; E = '  +3.00%  ' },
    @{ Row = 23; E = '  +3.11%  ' },
    @{ Row = 24; D = '268.84'; E = '  +2.41%  ' },
    @{ Row = 25; D = '2.35'; E = '  +7.95%  ' },
    @{ Row = 26; E = '  -0.50%  ' },
    @{ Row = 27; D = '7.62'; E = '  +9.59%  ' },
    @{ Row = 28; E = '  +5.69%  ' },
    @{ Row = 29; E = '  +0.24%  ' },
    @{ Row = 30; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '22.95'; E = '  +3.25%  ' },
    @{ Row = 31; B = 'InjectiveProtocol'; C = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D = '38.97'; E = '  +9.78%  ' },
    @{ Row = 32; D = '0.0936'; E = '  +9.16%  ' },
    @{ Row = 33; D = '169.93'; E = '  +2.97%  ' },
    @{ Row = 34; D = '2.99'; E = '  +16.27%  ' },
    @{ Row = 35; E = '  +2.27%  ' },
    @{ Row = 36; E = '  +4.65%  ' },
    @{ Row = 37; D = '4.83'; E = '  +7.71%  ' },
    @{ Row = 38; D = '3.08'; E = '  +13.34%  ' },
    @{ Row = 39; D = '0.0367'; E = '  +5.09%  ' },
    @{ Row = 40; D = '3.93'; E = '  +6.06%  ' },
    @{ Row = 41; E = '  +10.50%  ' },
    @{ Row = 42; D = '106.28'; E = '  +7.53%  ' },
    @{ Row = 43; D = '13.84'; E = '  +15.88%  ' },
    @{ Row = 44; E = '  +6.44%  ' },
    @{ Row = 45; D = '71.63'; E = '  +4.03%  ' },
    @{ Row = 46; E = '  -0.31%  ' },
    @{ Row = 47; D = '117.95'; E = '  +7.11%  ' },
    @{ Row = 48; D = '5.75'; E = '  +11.35%  ' },
    @{ Row = 49; B = 'MinaProtocolToken'; C = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'; D = '1.63'; E = '  +19.36%  ' },
    @{ Row = 50; B = 'ordi'; C = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'; D = '79.42'; E = '  +1.62%  ' },
    @{ Row = 51; D = '9.14'; E = '  +6.09%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B" + $u.Row).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $u.Row).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Range("D" + $u.Row).Value = "'" + $u.D }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $u.Row).Value = $u.E }
}
